# Automatische test-sync: 2025-06-30 19:39:50
# Adds a second test-mail row (row 3) to the "Logs" sheet and a matching
# summary row to the "Dashboard" sheet, then widens the conditional
# formatting ranges and the chart series ranges so they include row 3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Logs" sheet - append row 3 with the second test mail
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A3").Value = "Kun je deze order vandaag nog verwerken?"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("C3").Value = "Testmail #2: Kun je deze order vandaag nog verwerken?"
$logs.Range("D3").Value = "Bestelling / Levering"
$logs.Range("E3").Value = "Beste klant,`nDank u voor uw e-mail. Om uw verzoek om de order vandaag nog te verwerken te kunnen verwerken, hebben we meer informatie nodig. Kunt u ons het ordernummer en de specifieke items die u wilt bestellen doorgeven? Op die manier kunnen we uw verzoek zo snel mogelijk in behandeling nemen.`nMet vriendelijke groet,`n[Naam]  `nKlantenservice van [Bedrijfsnaam]"
$logs.Range("F3").Value = "2025-06-30 19:39:37"
$logs.Range("G3").Value = "Ja"
$logs.Range("H3").Value = "Ja"
$logs.Range("I3").Value = "Nee"
$logs.Range("J3").Value = "Nee"

# Widen the conditional-formatting ranges so every rule that used to
# apply to row 2 only now also covers row 3.
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $logs.Range($col + "2")
    $newRange = $logs.Range($col + "2:" + $col + "3")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 2) "Dashboard" sheet - append row 3 with the updated tally
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Bestelling / Levering"
$dash.Range("B3").Value = 1

# ---------------------------------------------------------------------
# 3) Chart on the "Dashboard" sheet - extend category/value references
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.XValues = "=Dashboard!`$A`$2:`$A`$3"
$series.Values = "=Dashboard!`$B`$2:`$B`$3"
